$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1443.8334
$ws.Range("I28").Value = 732.6
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 732.6
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -247.6
$ws.Range("N28").Value = -5970
$ws.Range("H64").Value = 3758.966
$ws.Range("I64").Value = 3548.276
$ws.Range("J64").Value = 3962.6333
$ws.Range("K64").Value = 3548.276
$ws.Range("L64").Value = 3962.6333
$ws.Range("M64").Value = -3300.276
$ws.Range("N64").Value = -4458.6333
$ws.Range("H67").Value = 3758.966
$ws.Range("I67").Value = 3548.276
$ws.Range("J67").Value = 3962.6333
$ws.Range("K67").Value = 3548.276
$ws.Range("L67").Value = 3962.6333
$ws.Range("M67").Value = -2690.276
$ws.Range("N67").Value = -5678.6333
$ws.Range("H87").Value = 90000
$ws.Range("J87").Value = 90000
$ws.Range("L87").Value = 90000
$ws.Range("N87").Value = -92496
$ws.Range("H90").Value = 90000
$ws.Range("J90").Value = 90000
$ws.Range("L90").Value = 270000
$ws.Range("N90").Value = -282480
$ws.Range("H111").Value = 553
$ws.Range("I111").Value = 553
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1659
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 1408
$ws.Range("N111").Value = $null
$ws.Range("H113").Value = 2125.25
$ws.Range("I113").Value = 2001.3334
$ws.Range("J113").Value = 2153.8462
$ws.Range("K113").Value = 2001.3334
$ws.Range("L113").Value = 2153.8462
$ws.Range("M113").Value = 1252.6666
$ws.Range("N113").Value = -8661.8462
$ws.Range("H116").Value = 1926.6666
$ws.Range("I116").Value = 1990
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 1990
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = 1452
$ws.Range("N116").Value = -8684
$ws.Range("H129").Value = 708.64
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 716.9796
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 2150.9388
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -12150.9388
$ws.Range("H132").Value = 2520.05
$ws.Range("I132").Value = 2079.302
$ws.Range("J132").Value = 5857.143
$ws.Range("K132").Value = 6237.906000000001
$ws.Range("L132").Value = 17571.429
$ws.Range("M132").Value = -3707.906000000001
$ws.Range("N132").Value = -22631.429
$ws.Range("H139").Value = 44274.285
$ws.Range("J139").Value = 44274.285
$ws.Range("L139").Value = 44274.285
$ws.Range("N139").Value = -54554.285
$ws.Range("H141").Value = 1296.25
$ws.Range("I141").Value = 1296.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3888.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1291.25
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1499.1154
$ws.Range("I2").Value = 1313.8948
$ws.Range("J2").Value = 2001.8572
$ws.Range("K2").Value = 1313.8948
$ws.Range("L2").Value = 2001.8572
$ws.Range("M2").Value = -1200.8948
$ws.Range("N2").Value = -2227.8572
$ws.Range("H32").Value = 15167042
$ws.Range("I32").Value = 21744386
$ws.Range("K32").Value = 21744386
$ws.Range("M32").Value = -21744099
$ws.Range("H116").Value = 1499.1154
$ws.Range("I116").Value = 1313.8948
$ws.Range("J116").Value = 2001.8572
$ws.Range("K116").Value = 1313.8948
$ws.Range("L116").Value = 2001.8572
$ws.Range("M116").Value = 980.1052
$ws.Range("N116").Value = -6589.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1499.1154
$ws.Range("I3").Value = 1313.8948
$ws.Range("J3").Value = 2001.8572
$ws.Range("K3").Value = 1313.8948
$ws.Range("L3").Value = 2001.8572
$ws.Range("M3").Value = -1199.8948
$ws.Range("N3").Value = -2229.8572
$ws.Range("H23").Value = 4750
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("H31").Value = 31285.715
$ws.Range("J31").Value = 31285.715
$ws.Range("L31").Value = 31285.715
$ws.Range("N31").Value = -31789.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500.0908
$ws.Range("I16").Value = 1103.6666
$ws.Range("J16").Value = 3023.75
$ws.Range("K16").Value = 1103.6666
$ws.Range("L16").Value = 3023.75
$ws.Range("M16").Value = -816.6666
$ws.Range("N16").Value = -3597.75
$ws.Range("H113").Value = 2500.0908
$ws.Range("I113").Value = 1103.6666
$ws.Range("J113").Value = 3023.75
$ws.Range("K113").Value = 1103.6666
$ws.Range("L113").Value = 3023.75
$ws.Range("M113").Value = 1066.3334
$ws.Range("N113").Value = -7363.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 802322.2
$ws.Range("I113").Value = 3135156.2
$ws.Range("J113").Value = 410.40625
$ws.Range("K113").Value = 9405468.600000001
$ws.Range("L113").Value = 1231.21875
$ws.Range("M113").Value = -9403298.600000001
$ws.Range("N113").Value = -5571.21875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 13469.1875
$ws.Range("I113").Value = 1230.1666
$ws.Range("J113").Value = 20812.6
$ws.Range("K113").Value = 1230.1666
$ws.Range("L113").Value = 20812.6
$ws.Range("M113").Value = 939.8334
$ws.Range("N113").Value = -25152.6
$ws.Range("H122").Value = 5556758
$ws.Range("I122").Value = 10000932
$ws.Range("J122").Value = 1540.5
$ws.Range("K122").Value = 30002796
$ws.Range("L122").Value = 4621.5
$ws.Range("M122").Value = -30000346
$ws.Range("N122").Value = -9521.5
$ws.Range("H126").Value = 5265023.5
$ws.Range("I126").Value = 7144093
$ws.Range("K126").Value = 21432279
$ws.Range("M126").Value = -21429809

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 999.6667
$ws.Range("J22").Value = 1316
$ws.Range("L22").Value = 1316
$ws.Range("N22").Value = -1906
$ws.Range("H27").Value = 999.6667
$ws.Range("J27").Value = 1316
$ws.Range("L27").Value = 1316
$ws.Range("N27").Value = -1530
$ws.Range("H61").Value = 1519.1538
$ws.Range("I61").Value = 1481.125
$ws.Range("J61").Value = 1580
$ws.Range("K61").Value = 1481.125
$ws.Range("L61").Value = 1580
$ws.Range("M61").Value = -1279.125
$ws.Range("N61").Value = -1984
$ws.Range("H113").Value = 1519.1538
$ws.Range("I113").Value = 1481.125
$ws.Range("J113").Value = 1580
$ws.Range("K113").Value = 1481.125
$ws.Range("L113").Value = 1580
$ws.Range("M113").Value = 688.875
$ws.Range("N113").Value = -5920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 608.9091
$ws.Range("I107").Value = 542.5714
$ws.Range("J107").Value = 725
$ws.Range("K107").Value = 1627.7142
$ws.Range("L107").Value = 2175
$ws.Range("M107").Value = 292.2857999999999
$ws.Range("N107").Value = -6015
